$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.847.73'
$ws.Range('E2').Value = '  -0.26%  '
$ws.Range('D3').Value = '2.542.56'
$ws.Range('E3').Value = '  -0.31%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '305.11'
$ws.Range('E5').Value = '  +1.92%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.86'
$ws.Range('E6').Value = '  +5.06%  '
$ws.Range('E7').Value = '  +0.86%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('E9').Value = '  -0.63%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.87'
$ws.Range('E10').Value = '  +1.60%  '
$ws.Range('E11').Value = '  +1.59%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.76'
$ws.Range('E12').Value = '  +0.31%  '
$ws.Range('E13').Value = '  -0.78%  '
$ws.Range('D14').Value = '2.933.27'
$ws.Range('E14').Value = '  -0.21%  '
$ws.Range('D15').Value = '2.563.93'
$ws.Range('E15').Value = '  +0.28%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.19'
$ws.Range('E16').Value = '  +6.95%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.877'
$ws.Range('E17').Value = '  -0.54%  '
$ws.Range('D18').Value = '42.858.78'
$ws.Range('E18').Value = '  -0.31%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.18'
$ws.Range('E19').Value = '  +3.55%  '
$ws.Range('D20').Value = '0.0₃0990'
$ws.Range('E20').Value = '  +0.78%  '
$ws.Range('E21').Value = '  +0.18%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '71.65'
$ws.Range('E22').Value = '  -0.52%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '254.24'
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('E24').Value = '  +0.44%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.07'
$ws.Range('E25').Value = '  -2.68%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '27.73'
$ws.Range('E26').Value = '  -4.44%  '
$ws.Range('E27').Value = '  +0.12%  '
$ws.Range('E28').Value = '  +9.82%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.20'
$ws.Range('E29').Value = '  -0.52%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '38.86'
$ws.Range('E30').Value = '  +5.58%  '
$ws.Range('E31').Value = '  +0.77%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '157.66'
$ws.Range('E32').Value = '  +3.62%  '
$ws.Range('E33').Value = '  +0.08%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0803'
$ws.Range('E34').Value = '  +1.23%  '
$ws.Range('B35').Value = 'Celestia'
$ws.Range('C35').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '19.03'
$ws.Range('E35').Value = '  +7.58%  '
$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.29'
$ws.Range('E36').Value = '  -2.20%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.64'
$ws.Range('E37').Value = '  -4.25%  '
$ws.Range('E38').Value = '  +0.71%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '24.57'
$ws.Range('E39').Value = '  +6.79%  '
$ws.Range('E40').Value = '  +0.76%  '
$ws.Range('E41').Value = '  +6.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.46'
$ws.Range('E42').Value = '  +1.26%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.89'
$ws.Range('E43').Value = '  +0.51%  '
$ws.Range('D44').Value = '2.085.11'
$ws.Range('E44').Value = '  -1.26%  '
$ws.Range('E45').Value = '  -1.97%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.999'
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '86.34'
$ws.Range('E47').Value = '  +2.57%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.08'
$ws.Range('E48').Value = '  -0.14%  '
$ws.Range('D49').Value = '2.790.50'
$ws.Range('E49').Value = '  -0.19%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '73.61'
$ws.Range('E50').Value = '  +6.15%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.192'
$ws.Range('E51').Value = '  +1.62%  '
